$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the "Micro results" row robustly by scanning column 1 text
$targetRow = -1
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $labelCell = $t.Cell($r, 1)
    if ($labelCell.Range.Text -like "*Micro results*") {
        $targetRow = $r
        break
    }
}

if ($targetRow -eq -1) {
    throw "Could not find the 'Micro results' row"
}

$cell = $t.Cell($targetRow, 2)
$p = $cell.Range.Paragraphs.First

$xmlFrag = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
            <w:r/>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>27/08 – TB MICROSCOPY – **Negative**</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="preserve">  Summary: No mycobacteria observed.</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>27/08 – RESP. CULT AND MICRO – **Negative**</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="preserve">  Summary: Upper respiratory flora only.</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>27/08 – TB CULTURE – SPUTUM MYCOBACTERIAL CULTURE: NO GROWTH AFTER 8 WEEKS</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>26/08 – BLC – ;Other (specify site in Clinical Details) NO GROWTH AFTER 5 DAYS</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>26/08 – BLC – RED PORT NO GROWTH AFTER 5 DAYS</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>25/08 – COPIES/ML – Negative</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>25/08 – EBV DNA – Negative</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>25/08 – CMV DNA – Negative</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>24/08 – BLC – ;Other (specify site in Clinical Details) NO GROWTH AFTER 5 DAYS</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>24/08 – BLC – PERIPHERAL–RIGHT NO GROWTH AFTER 5 DAYS</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>24/08 – BLC – RED PORT NO GROWTH AFTER 5 DAYS</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>22/08 – BLC – PERIPHERAL–RIGHT NO GROWTH AFTER 5 DAYS</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>22/08 – BLC – RED PORT NO GROWTH AFTER 5 DAYS</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>22/08 – BLC – ;Pic Line NO GROWTH AFTER 5 DAYS</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>22/08 – RESPIRATORY PCR – **Negative**</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="preserve">  Summary: No microorganisms detected by PCR.</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>--------Previous result (1 year)--------</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>24/07 – EBV VCA IgG – Positive</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>19/06 – RESPIRATORY PCR – **Positive** + **Summary:** Human Rhinovirus/Enterovirus detected.</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>12/06 – EBV VCA IgG – Positive</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="preserve">20/02 – RESPIRATORY PCR – **Positive**  </w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>**Summary:** RSV detected by PCR.</w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="preserve">04/02 – RESPIRATORY PCR – **Positive**  </w:t>
            </w:r>
          </w:p>
<w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="0000FF"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t>**Summary:** RSV detected.</w:t>
            </w:r>
          </w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$p.Range.InsertXML($xmlFrag)

Write-Host "Micro results cell rebuilt."
